$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.600.96'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = '2.586.05'
$ws.Range('E3').Value = '  +10.41%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '305.02'
$ws.Range('E5').Value = '  +1.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.30'
$ws.Range('E6').Value = '  +3.06%  '
$ws.Range('E7').Value = '  +5.91%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.565'
$ws.Range('E9').Value = '  +10.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.46'
$ws.Range('E10').Value = '  +10.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0830'
$ws.Range('E11').Value = '  +4.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.98'
$ws.Range('E12').Value = '  +12.14%  '
$ws.Range('D13').Value = '2.984.78'
$ws.Range('E13').Value = '  +10.64%  '
$ws.Range('E14').Value = '  +2.44%  '
$ws.Range('D15').Value = '2.586.50'
$ws.Range('E15').Value = '  +10.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.897'
$ws.Range('E16').Value = '  +11.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '15.01'
$ws.Range('E17').Value = '  +9.59%  '
$ws.Range('D18').Value = '46.769.77'
$ws.Range('E18').Value = '  +1.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.91'
$ws.Range('E19').Value = '  +9.42%  '
$ws.Range('D20').Value = '0.0₃0998'
$ws.Range('E20').Value = '  +3.57%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.56'
$ws.Range('E21').Value = '  +9.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '69.83'
$ws.Range('E22').Value = '  +4.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '254.05'
$ws.Range('E23').Value = '  +3.76%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.97'
$ws.Range('E24').Value = '  +5.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.12'
$ws.Range('E25').Value = '  +11.44%  '
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '26.10'
$ws.Range('E27').Value = '  +23.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '40.79'
$ws.Range('E28').Value = '  +1.95%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.38'
$ws.Range('E29').Value = '  +7.46%  '
$ws.Range('E30').Value = '  +3.64%  '
$ws.Range('E31').Value = '  +3.74%  '
$ws.Range('E32').Value = '  +4.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.93'
$ws.Range('E33').Value = '  +8.23%  '
$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('C34').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.22'
$ws.Range('E34').Value = '  +21.19%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0834'
$ws.Range('E35').Value = '  +7.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '148.66'
$ws.Range('E36').Value = '  +3.12%  '
$ws.Range('E37').Value = '  +8.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.120'
$ws.Range('E38').Value = '  +3.60%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '16.26'
$ws.Range('E39').Value = '  +8.80%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.22'
$ws.Range('E40').Value = '  +9.00%  '
$ws.Range('B41').Value = 'NEARProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.58'
$ws.Range('E41').Value = '  +11.22%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0326'
$ws.Range('E42').Value = '  +8.59%  '
$ws.Range('D43').Value = '2.034.75'
$ws.Range('E43').Value = '  +10.43%  '
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.56'
$ws.Range('E45').Value = '  +2.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.56'
$ws.Range('E46').Value = '  +33.99%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.84'
$ws.Range('E47').Value = '  +0.25%  '
$ws.Range('E48').Value = '  +7.42%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.841.62'
$ws.Range('E49').Value = '  +10.54%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '106.89'
$ws.Range('E50').Value = '  +10.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.78'
$ws.Range('E51').Value = '  +9.31%  '
